# "removed default client from excel file"
#
# The Client-Template workbook ships with a sample/default client
# ("Goood Media") pre-filled in row 2 of the "Client Info" sheet, plus
# hyperlinks pointing at that sample client's website/email. Strip that
# seed data out so the template starts blank, and leave the "Client Info"
# sheet as the active tab/selection (instead of "Domain & Hosting").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Client Info")

# Row 2 holds the default client's sample values (Company/Name/Email/
# Phone/Address/Role/...). Blow away both the values and the formatting
# that was only there to show off the sample row.
$ws1.Range("A2:B2").Clear()
$ws1.Range("C2:F2").ClearContents()
$ws1.Range("G2:H2").Clear()

# The sample company's website (C2) and email (E2) were hyperlinked —
# drop those links along with the data.
$ws1.Hyperlinks.Delete()

# Bring "Client Info" to the front (it was "Domain & Hosting" before)
# with the now-empty client row selected, ready for real data entry.
$ws1.Activate()
$ws1.Rows("2:2").Select()
